$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$uChar1_2083 = [char]0x2083

function Set-TextValue {
    param($Sheet, [string]$Addr, [string]$Val)
    $cell = $Sheet.Range($Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.ClearFormats()
}

Set-TextValue $ws "D2" "26.486.23"
Set-TextValue $ws "E2" "  -0.13%  "
Set-TextValue $ws "D3" "1.625.56"
Set-TextValue $ws "E3" "  +0.35%  "
Set-TextValue $ws "D4" "1.00"
Set-TextValue $ws "E4" "  +0.15%  "
Set-TextValue $ws "D5" "213.59"
Set-TextValue $ws "E5" "  -0.33%  "
Set-TextValue $ws "E6" "  -0.66%  "
Set-TextValue $ws "E7" "  +0.16%  "
Set-TextValue $ws "E8" "  -0.09%  "
Set-TextValue $ws "E9" "  -0.17%  "
Set-TextValue $ws "D10" "19.20"
Set-TextValue $ws "D11" "0.0854"
Set-TextValue $ws "E11" "  -0.28%  "
Set-TextValue $ws "D12" "1.854.37"
Set-TextValue $ws "E12" "  +0.39%  "
Set-TextValue $ws "D13" "1.624.90"
Set-TextValue $ws "E13" "  +0.35%  "
Set-TextValue $ws "E14" "  -0.09%  "
Set-TextValue $ws "E15" "  -0.41%  "
Set-TextValue $ws "D16" "63.98"
Set-TextValue $ws "E16" "  -1.31%  "
Set-TextValue $ws "D17" "234.54"
Set-TextValue $ws "E17" "  +0.88%  "
Set-TextValue $ws "D18" "26.497.00"
Set-TextValue $ws "E18" "  -0.03%  "
Set-TextValue $ws "D19" "7.76"
Set-TextValue $ws "E19" "  +2.39%  "
Set-TextValue $ws "D20" "0.0${uChar1_2083}0726"
Set-TextValue $ws "E20" "  -0.15%  "
Set-TextValue $ws "D21" "0.999"
Set-TextValue $ws "E21" "  +0.09%  "
Set-TextValue $ws "E22" "  -1.53%  "
Set-TextValue $ws "E23" "  +2.75%  "
Set-TextValue $ws "D24" "9.13"
Set-TextValue $ws "E24" "  -0.07%  "
Set-TextValue $ws "D25" "146.98"
Set-TextValue $ws "E25" "  +1.08%  "
Set-TextValue $ws "E26" "  +0.08%  "
Set-TextValue $ws "E27" "  +0.32%  "
Set-TextValue $ws "D28" "0.113"
Set-TextValue $ws "E28" "  -0.56%  "
Set-TextValue $ws "D29" "15.64"
Set-TextValue $ws "E29" "  +0.59%  "
Set-TextValue $ws "D30" "0.0497"
Set-TextValue $ws "E30" "  -0.21%  "
Set-TextValue $ws "E31" "  -0.25%  "
Set-TextValue $ws "D32" "1.520.81"
Set-TextValue $ws "E32" "  +4.82%  "
Set-TextValue $ws "E33" "  +0.87%  "
Set-TextValue $ws "E34" "  -0.79%  "
Set-TextValue $ws "E35" "  +2.75%  "
Set-TextValue $ws "E36" "  +0.22%  "
Set-TextValue $ws "D37" "0.567"
Set-TextValue $ws "E37" "  +1.16%  "
Set-TextValue $ws "D38" "0.0167"
Set-TextValue $ws "E38" "  -0.33%  "
Set-TextValue $ws "E39" "  -0.51%  "
Set-TextValue $ws "E40" "  -0.08%  "
Set-TextValue $ws "E41" "  +0.06%  "
Set-TextValue $ws "E42" "  +0.34%  "
Set-TextValue $ws "D43" "1.765.16"
Set-TextValue $ws "E43" "  +0.47%  "
Set-TextValue $ws "D44" "62.97"
Set-TextValue $ws "E44" "  +1.11%  "
Set-TextValue $ws "E45" "  -0.16%  "
Set-TextValue $ws "D46" "0.906"
Set-TextValue $ws "E46" "  -1.21%  "
Set-TextValue $ws "D47" "90.06"
Set-TextValue $ws "E47" "  +1.96%  "
Set-TextValue $ws "E48" "  +0.41%  "
Set-TextValue $ws "B49" "Cronos"
Set-TextValue $ws "C49" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D49" "0.0501"
Set-TextValue $ws "E49" "  -0.36%  "
Set-TextValue $ws "B50" "Algorand"
Set-TextValue $ws "C50" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws "D50" "0.0964"
Set-TextValue $ws "E50" "  -0.09%  "
Set-TextValue $ws "B51" "EnergySwap"
Set-TextValue $ws "C51" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D51" "7.50"
Set-TextValue $ws "E51" "  -0.16%  "
